$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 10 (old MuSCs-sending rows removed)
$ws.Range("A10:T10").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.068624
$ws.Range("H2").Value = 0.205872
$ws.Range("I2").Value = 0.01198115042951486
$ws.Range("J2").Value = 0.01198115042951486
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 19.524618
$ws.Range("N2").Value = 58.573854
$ws.Range("O2").Value = 0.4154885426712971
$ws.Range("P2").Value = 0.4539723485554654
$ws.Range("Q2").Value = 1.339857385632
$ws.Range("R2").Value = 12.058716470688
$ws.Range("S2").Value = 0.004978030731484715
$ws.Range("T2").Value = 0.005439110998883185

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.068624
$ws.Range("H3").Value = 0.205872
$ws.Range("I3").Value = 0.01198115042951486
$ws.Range("J3").Value = 0.01198115042951486
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.24435933333334
$ws.Range("N3").Value = 45.73307800000001
$ws.Range("O3").Value = 0.324403614112412
$ws.Range("P3").Value = 0.3544508583357054
$ws.Range("Q3").Value = 1.046128914890667
$ws.Range("R3").Value = 9.415160234016001
$ws.Range("S3").Value = 0.003886728500559098
$ws.Range("T3").Value = 0.004246729053590748

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.068624
$ws.Range("H4").Value = 0.205872
$ws.Range("I4").Value = 0.01198115042951486
$ws.Range("J4").Value = 0.01198115042951486
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 11.9507005
$ws.Range("N4").Value = 23.901401
$ws.Range("O4").Value = 0.2543137660693869
$ws.Range("P4").Value = 0.1852460510065796
$ws.Range("Q4").Value = 0.820104871112
$ws.Range("R4").Value = 4.920629226672
$ws.Range("S4").Value = 0.003046971487573777
$ws.Range("T4").Value = 0.002219460803583414

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.068624
$ws.Range("H5").Value = 0.205872
$ws.Range("I5").Value = 0.01198115042951486
$ws.Range("J5").Value = 0.01198115042951486
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.272275
$ws.Range("N5").Value = 0.816825
$ws.Range("O5").Value = 0.005794077146903843
$ws.Range("P5").Value = 0.006330742102249548
$ws.Range("Q5").Value = 0.0186845996
$ws.Range("R5").Value = 0.1681613964
$ws.Range("S5").Value = 0.00006941970989726922
$ws.Range("T5").Value = 0.00007584957345751498

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.659039666666668
$ws.Range("H6").Value = 16.977119
$ws.Range("I6").Value = 0.9880188495704851
$ws.Range("J6").Value = 0.9880188495704851
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 19.524618
$ws.Range("N6").Value = 58.573854
$ws.Range("O6").Value = 0.4154885426712971
$ws.Range("P6").Value = 0.4539723485554654
$ws.Range("Q6").Value = 110.490587738514
$ws.Range("R6").Value = 994.4152896466261
$ws.Range("S6").Value = 0.4105105119398124
$ws.Range("T6").Value = 0.4485332375565823

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.659039666666668
$ws.Range("H7").Value = 16.977119
$ws.Range("I7").Value = 0.9880188495704851
$ws.Range("J7").Value = 0.9880188495704851
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 15.24435933333334
$ws.Range("N7").Value = 45.73307800000001
$ws.Range("O7").Value = 0.324403614112412
$ws.Range("P7").Value = 0.3544508583357054
$ws.Range("Q7").Value = 86.26843416025358
$ws.Range("R7").Value = 776.4159074422822
$ws.Range("S7").Value = 0.3205168856118529
$ws.Range("T7").Value = 0.3502041292821146

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.659039666666668
$ws.Range("H8").Value = 16.977119
$ws.Range("I8").Value = 0.9880188495704851
$ws.Range("J8").Value = 0.9880188495704851
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 11.9507005
$ws.Range("N8").Value = 23.901401
$ws.Range("O8").Value = 0.2543137660693869
$ws.Range("P8").Value = 0.1852460510065796
$ws.Range("Q8").Value = 67.62948817395318
$ws.Range("R8").Value = 405.7769290437191
$ws.Range("S8").Value = 0.2512667945818132
$ws.Range("T8").Value = 0.1830265902029962

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.659039666666668
$ws.Range("H9").Value = 16.977119
$ws.Range("I9").Value = 0.9880188495704851
$ws.Range("J9").Value = 0.9880188495704851
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.272275
$ws.Range("N9").Value = 0.816825
$ws.Range("O9").Value = 0.005794077146903843
$ws.Range("P9").Value = 0.006330742102249548
$ws.Range("Q9").Value = 1.540815025241667
$ws.Range("R9").Value = 13.867335227175
$ws.Range("S9").Value = 0.005724657437006574
$ws.Range("T9").Value = 0.006254892528792032
